$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 109.1118546666667
$ws.Range("H2").Value = 327.335564
$ws.Range("I2").Value = 0.3029068882986101
$ws.Range("J2").Value = 0.3029068882986101
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 14.81590833333333
$ws.Range("N2").Value = 44.447725
$ws.Range("O2").Value = 0.2626950900477981
$ws.Range("P2").Value = 0.2626950900477981
$ws.Range("Q2").Value = 1616.591236821322
$ws.Range("R2").Value = 14549.3211313919
$ws.Range("S2").Value = 0.07957215229770168
$ws.Range("T2").Value = 0.07957215229770168

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 109.1118546666667
$ws.Range("H3").Value = 327.335564
$ws.Range("I3").Value = 0.3029068882986101
$ws.Range("J3").Value = 0.3029068882986101
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 34.77831566666666
$ws.Range("N3").Value = 104.334947
$ws.Range("O3").Value = 0.6166407458941315
$ws.Range("P3").Value = 0.6166407458941315
$ws.Range("Q3").Value = 3794.72652457279
$ws.Range("R3").Value = 34152.53872115511
$ws.Range("S3").Value = 0.1867847295369253
$ws.Range("T3").Value = 0.1867847295369253

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 109.1118546666667
$ws.Range("H4").Value = 327.335564
$ws.Range("I4").Value = 0.3029068882986101
$ws.Range("J4").Value = 0.3029068882986101
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.805415333333333
$ws.Range("N4").Value = 20.416246
$ws.Range("O4").Value = 0.1206641640580704
$ws.Range("P4").Value = 0.1206641640580704
$ws.Range("Q4").Value = 742.5514887969716
$ws.Range("R4").Value = 6682.963399172745
$ws.Range("S4").Value = 0.03655000646398309
$ws.Range("T4").Value = 0.03655000646398309

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 227.6338753333333
$ws.Range("H5").Value = 682.901626
$ws.Range("I5").Value = 0.6319374650831437
$ws.Range("J5").Value = 0.6319374650831437
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 14.81590833333333
$ws.Range("N5").Value = 44.447725
$ws.Range("O5").Value = 0.2626950900477981
$ws.Range("P5").Value = 0.2626950900477981
$ws.Range("Q5").Value = 3372.602630500094
$ws.Range("R5").Value = 30353.42367450085
$ws.Range("S5").Value = 0.1660068692945937
$ws.Range("T5").Value = 0.1660068692945937

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 227.6338753333333
$ws.Range("H6").Value = 682.901626
$ws.Range("I6").Value = 0.6319374650831437
$ws.Range("J6").Value = 0.6319374650831437
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 34.77831566666666
$ws.Range("N6").Value = 104.334947
$ws.Range("O6").Value = 0.6166407458941315
$ws.Range("P6").Value = 0.6166407458941315
$ws.Range("Q6").Value = 7916.722772769313
$ws.Range("R6").Value = 71250.50495492382
$ws.Range("S6").Value = 0.3896783898273165
$ws.Range("T6").Value = 0.3896783898273165

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 227.6338753333333
$ws.Range("H7").Value = 682.901626
$ws.Range("I7").Value = 0.6319374650831437
$ws.Range("J7").Value = 0.6319374650831437
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 6.805415333333333
$ws.Range("N7").Value = 20.416246
$ws.Range("O7").Value = 0.1206641640580704
$ws.Range("P7").Value = 0.1206641640580704
$ws.Range("Q7").Value = 1549.143065579555
$ws.Range("R7").Value = 13942.287590216
$ws.Range("S7").Value = 0.0762522059612336
$ws.Range("T7").Value = 0.0762522059612336

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.194568
$ws.Range("H8").Value = 0.583704
$ws.Range("I8").Value = 0.0005401428435299865
$ws.Range("J8").Value = 0.0005401428435299865
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 14.81590833333333
$ws.Range("N8").Value = 44.447725
$ws.Range("O8").Value = 0.2626950900477981
$ws.Range("P8").Value = 0.2626950900477981
$ws.Range("Q8").Value = 2.8827016526
$ws.Range("R8").Value = 25.9443148734
$ws.Range("S8").Value = 0.0001418928729197835
$ws.Range("T8").Value = 0.0001418928729197835

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.194568
$ws.Range("H9").Value = 0.583704
$ws.Range("I9").Value = 0.0005401428435299865
$ws.Range("J9").Value = 0.0005401428435299865
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 34.77831566666666
$ws.Range("N9").Value = 104.334947
$ws.Range("O9").Value = 0.6166407458941315
$ws.Range("P9").Value = 0.6166407458941315
$ws.Range("Q9").Value = 6.766747322631999
$ws.Range("R9").Value = 60.900725903688
$ws.Range("S9").Value = 0.0003330740859237081
$ws.Range("T9").Value = 0.0003330740859237081

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.194568
$ws.Range("H10").Value = 0.583704
$ws.Range("I10").Value = 0.0005401428435299865
$ws.Range("J10").Value = 0.0005401428435299865
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 6.805415333333333
$ws.Range("N10").Value = 20.416246
$ws.Range("O10").Value = 0.1206641640580704
$ws.Range("P10").Value = 0.1206641640580704
$ws.Range("Q10").Value = 1.324116050576
$ws.Range("R10").Value = 11.917044455184
$ws.Range("S10").Value = 0.00006517588468649495
$ws.Range("T10").Value = 0.00006517588468649495

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.4517016666666667
$ws.Range("H11").Value = 1.355105
$ws.Range("I11").Value = 0.001253975076377243
$ws.Range("J11").Value = 0.001253975076377243
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 14.81590833333333
$ws.Range("N11").Value = 44.447725
$ws.Range("O11").Value = 0.2626950900477981
$ws.Range("P11").Value = 0.2626950900477981
$ws.Range("Q11").Value = 6.692370487347222
$ws.Range("R11").Value = 60.231334386125
$ws.Range("S11").Value = 0.0003294130956066143
$ws.Range("T11").Value = 0.0003294130956066143

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.4517016666666667
$ws.Range("H12").Value = 1.355105
$ws.Range("I12").Value = 0.001253975076377243
$ws.Range("J12").Value = 0.001253975076377243
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 34.77831566666666
$ws.Range("N12").Value = 104.334947
$ws.Range("O12").Value = 0.6166407458941315
$ws.Range("P12").Value = 0.6166407458941315
$ws.Range("Q12").Value = 15.70942315049278
$ws.Range("R12").Value = 141.384808354435
$ws.Range("S12").Value = 0.0007732521264299137
$ws.Range("T12").Value = 0.0007732521264299137

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.4517016666666667
$ws.Range("H13").Value = 1.355105
$ws.Range("I13").Value = 0.001253975076377243
$ws.Range("J13").Value = 0.001253975076377243
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 6.805415333333333
$ws.Range("N13").Value = 20.416246
$ws.Range("O13").Value = 0.1206641640580704
$ws.Range("P13").Value = 0.1206641640580704
$ws.Range("Q13").Value = 3.074017448425556
$ws.Range("R13").Value = 27.66615703583
$ws.Range("S13").Value = 0.000151309854340715
$ws.Range("T13").Value = 0.000151309854340715

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.7729569999999999
$ws.Range("H14").Value = 2.318871
$ws.Range("I14").Value = 0.002145816331084288
$ws.Range("J14").Value = 0.002145816331084288
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 14.81590833333333
$ws.Range("N14").Value = 44.447725
$ws.Range("O14").Value = 0.2626950900477981
$ws.Range("P14").Value = 0.2626950900477981
$ws.Range("Q14").Value = 11.45206005760833
$ws.Range("R14").Value = 103.068540518475
$ws.Range("S14").Value = 0.0005636954143202227
$ws.Range("T14").Value = 0.0005636954143202227

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.7729569999999999
$ws.Range("H15").Value = 2.318871
$ws.Range("I15").Value = 0.002145816331084288
$ws.Range("J15").Value = 0.002145816331084288
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 34.77831566666666
$ws.Range("N15").Value = 104.334947
$ws.Range("O15").Value = 0.6166407458941315
$ws.Range("P15").Value = 0.6166407458941315
$ws.Range("Q15").Value = 26.88214254275966
$ws.Range("R15").Value = 241.939282884837
$ws.Range("S15").Value = 0.001323197782951624
$ws.Range("T15").Value = 0.001323197782951624

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.7729569999999999
$ws.Range("H16").Value = 2.318871
$ws.Range("I16").Value = 0.002145816331084288
$ws.Range("J16").Value = 0.002145816331084288
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.805415333333333
$ws.Range("N16").Value = 20.416246
$ws.Range("O16").Value = 0.1206641640580704
$ws.Range("P16").Value = 0.1206641640580704
$ws.Range("Q16").Value = 5.260293419807333
$ws.Range("R16").Value = 47.342640778266
$ws.Range("S16").Value = 0.0002589231338124413
$ws.Range("T16").Value = 0.0002589231338124413

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 22.05086833333333
$ws.Range("H17").Value = 66.15260499999999
$ws.Range("I17").Value = 0.06121571236725463
$ws.Range("J17").Value = 0.06121571236725463
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 14.81590833333333
$ws.Range("N17").Value = 44.447725
$ws.Range("O17").Value = 0.2626950900477981
$ws.Range("P17").Value = 0.2626950900477981
$ws.Range("Q17").Value = 326.7036438970694
$ws.Range("R17").Value = 2940.332795073625
$ws.Range("S17").Value = 0.01608106707265606
$ws.Range("T17").Value = 0.01608106707265606

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 22.05086833333333
$ws.Range("H18").Value = 66.15260499999999
$ws.Range("I18").Value = 0.06121571236725463
$ws.Range("J18").Value = 0.06121571236725463
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 34.77831566666666
$ws.Range("N18").Value = 104.334947
$ws.Range("O18").Value = 0.6166407458941315
$ws.Range("P18").Value = 0.6166407458941315
$ws.Range("Q18").Value = 766.8920596207704
$ws.Range("R18").Value = 6902.028536586934
$ws.Range("S18").Value = 0.03774810253458451
$ws.Range("T18").Value = 0.03774810253458451

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 22.05086833333333
$ws.Range("H19").Value = 66.15260499999999
$ws.Range("I19").Value = 0.06121571236725463
$ws.Range("J19").Value = 0.06121571236725463
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 6.805415333333333
$ws.Range("N19").Value = 20.416246
$ws.Range("O19").Value = 0.1206641640580704
$ws.Range("P19").Value = 0.1206641640580704
$ws.Range("Q19").Value = 150.0653174689811
$ws.Range("R19").Value = 1350.58785722083
$ws.Range("S19").Value = 0.007386542760014062
$ws.Range("T19").Value = 0.007386542760014062
